# Add a new "M_PL" (profit) column to the income table, inserted right
# after the existing "M_ETR" column (i.e. before the old column C).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at C; this shifts the old C:J data to D:K
# and keeps the header cell's style/border formatting intact.
$ws.Columns("C:C").Insert()

# New header for the inserted column.
$ws.Range("C1").Value = "M_PL"

# New profit figures for the inserted column (one per group row).
$ws.Range("C2").Value = 1008209699708
$ws.Range("C3").Value = 3140810
$ws.Range("C4").Value = 21277927825
$ws.Range("C5").Value = 353160988340
$ws.Range("C6").Value = 49527932043
